$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Arduino interpretation" column (D) for rows 2-5.
# Write order matters for shared-string allocation order (matches the
# target workbook's shared-strings table, which appends new unique
# strings in this sequence: row3, row4, row5, row2).
$ws.Range("D3").Value = "Karte an diesem Wochentag ungültig"
$ws.Range("D4").Value = "Ihre Karten wurde innerhalb 3 Stunden bereits verwendet"
$ws.Range("D5").Value = "Ihre Karten ID ist nicht bekannt"
$ws.Range("D2").Value = "Karte abglaufen"

# Update the saved selection/active cell to D5, matching the workbook
# view state captured at save time.
$ws.Range("D5").Select()
